$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G42").Value = "Fallo"
$ws.Range("H42").Value = -1

$ws.Range("G48").Value = "Acierto"
$ws.Range("H48").Value = 1.62

$ws.Range("G49").Value = "Acierto"
$ws.Range("H49").Value = 0.83

$ws.Range("G50").Value = "Acierto"
$ws.Range("H50").Value = 0.83

$ws.Range("G51").Value = "Acierto"
$ws.Range("H51").Value = 1.75

$ws.Range("G52").Value = "Fallo"
$ws.Range("H52").Value = -1
